$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "deuteron" with "d" in the target column (G2:G11)
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "deuteron") {
        $cell.Value = "d"
    }
}

# Make the header row bold
$ws.Range("A1:K1").Font.Bold = $true

# Update the view selection to match the edited file
$ws.Range("G19").Select()
